$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5416.988499999999
$ws.Range("B3").Value = 5210.142
$ws.Range("B4").Value = 5131.699999999999
$ws.Range("B5").Value = 5049.197999999999
$ws.Range("B6").Value = 5084.9925
$ws.Range("B7").Value = 5380.1055
$ws.Range("B8").Value = 6413.1515
$ws.Range("B9").Value = 7207.8405
$ws.Range("B10").Value = 9489.200000000001
$ws.Range("B11").Value = 14318.1885
$ws.Range("B12").Value = 15340.192
$ws.Range("B13").Value = 15589.8715
$ws.Range("B14").Value = 15359.3685
$ws.Range("B15").Value = 15511.44
$ws.Range("B16").Value = 15666.4375
$ws.Range("B17").Value = 16125.417
$ws.Range("B18").Value = 16406.873
$ws.Range("B19").Value = 16163.1645
$ws.Range("B20").Value = 15659.7175
$ws.Range("B21").Value = 13764.5865
$ws.Range("B22").Value = 11506.208
$ws.Range("B23").Value = 8102.793999999999
$ws.Range("B24").Value = 6308.112999999999
$ws.Range("B25").Value = 5711.482
